$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - updated TPM-derived values
$ws.Range("M2").Value = 7.374758333333332
$ws.Range("N2").Value = 22.124275
$ws.Range("O2").Value = 0.2902113952021427
$ws.Range("P2").Value = 0.2902113952021427
$ws.Range("Q2").Value = 0.6432165888222221
$ws.Range("R2").Value = 5.7889492994
$ws.Range("S2").Value = 0.2902113952021427
$ws.Range("T2").Value = 0.2902113952021427

# Row 3 - updated specificity values (M/N unchanged)
$ws.Range("O3").Value = 0.07140925709849509
$ws.Range("P3").Value = 0.07140925709849508
$ws.Range("S3").Value = 0.07140925709849509
$ws.Range("T3").Value = 0.07140925709849508

# Row 4 - updated specificity values (M/N unchanged)
$ws.Range("O4").Value = 0.6383793476993622
$ws.Range("P4").Value = 0.6383793476993621
$ws.Range("S4").Value = 0.6383793476993622
$ws.Range("T4").Value = 0.6383793476993621
